$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 76

# Column A holds a date-like label ("01-07-2021"). A plain .Value
# assignment gets auto-converted by Excel into a date serial number, which
# isn't what the source data has (it's the literal text used throughout
# the rest of column A). Force it in as text: put it in via a formula
# that evaluates to that exact string, then Copy/PasteSpecial values-only
# over itself so the cell ends up holding a plain text value (no formula,
# no leftover number-format/style changes).
$ws.Cells.Item($row, 1).Formula = "=""01-07-2021"""
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 2).Value = 27333
$ws.Cells.Item($row, 3).Value = 10020
$ws.Cells.Item($row, 4).Value = 1410
$ws.Cells.Item($row, 5).Value = 5925
$ws.Cells.Item($row, 6).Value = 2684
$ws.Cells.Item($row, 7).Value = 17314
$ws.Cells.Item($row, 8).Value = 11857
$ws.Cells.Item($row, 9).Value = 5456
